$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.491.83"
$ws.Range("E2").Value = "  -3.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.754.88"
$ws.Range("E3").Value = "  -2.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.71"
$ws.Range("E5").Value = "  -1.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4405"
$ws.Range("E7").Value = "  -1.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3701"
$ws.Range("E8").Value = "  -2.01%  "

$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07670"
$ws.Range("E10").Value = "  +2.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.113"
$ws.Range("E11").Value = "  -3.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.58"
$ws.Range("E13").Value = "  -4.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.149"
$ws.Range("E14").Value = "  -2.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.411"
$ws.Range("E15").Value = "  -2.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.759.83"
$ws.Range("E16").Value = "  -2.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.40"
$ws.Range("E17").Value = "  +12.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001073"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06250"
$ws.Range("E19").Value = "  -8.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.38"
$ws.Range("E21").Value = "  -1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.170"
$ws.Range("E22").Value = "  -2.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5305"
$ws.Range("E23").Value = "  -2.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.542.27"
$ws.Range("E24").Value = "  -3.08%  "

$ws.Range("E25").Value = "  -2.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.299"
$ws.Range("E26").Value = "  -4.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.51"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.03"
$ws.Range("E28").Value = "  -0.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.295"
$ws.Range("E29").Value = "  -2.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.957.93"
$ws.Range("E30").Value = "  -2.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.42"
$ws.Range("E31").Value = "  -3.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.175"
$ws.Range("E32").Value = "  -6.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.713"
$ws.Range("E33").Value = "  -1.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09174"
$ws.Range("E34").Value = "  -1.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.644"
$ws.Range("E35").Value = "  -8.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.61"
$ws.Range("E36").Value = "  +3.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02313"
$ws.Range("E37").Value = "  -1.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2156"
$ws.Range("E38").Value = "  -5.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06109"
$ws.Range("E39").Value = "  -4.09%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.052"
$ws.Range("E40").Value = "  -1.95%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6433"
$ws.Range("E41").Value = "  -2.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.175"
$ws.Range("E42").Value = "  -2.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.937"
$ws.Range("E43").Value = "  -2.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.390"
$ws.Range("E45").Value = "  -4.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.66"
$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5961"
$ws.Range("E47").Value = "  -2.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.717"
$ws.Range("E48").Value = "  -2.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.16"
$ws.Range("E49").Value = "  -1.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.980"
$ws.Range("E50").Value = "  -2.52%  "

$ws.Range("E51").Value = "  -3.06%  "
